# Fruta / hortaliza, semanal
# Update weekly price records for "Feria Lagunitas de Puerto Montt - Kiwi" and
# append the newest batch of rows (119-122), pushing the previously-last
# observations down, mirroring the upstream weekly refresh of this dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, Fecha(D), Calidad(L), Volumen(M), PrecioMin(N),
#             PrecioMax(O), PrecioPromedio(P), UnidadComercializacion(Q),
#             Origen(R), PrecioKg(S), KgUnidad(T)
$rows = @(
    @(101, 44474, 'Especial', 300, 21000, 21000, 21000, '$/caja 15 kilos', 'Provincia de Curicó', 1400, 15),
    @(102, 44474, 'Primera',  200, 15000, 15000, 15000, '$/caja 15 kilos', 'Provincia de Curicó', 1000, 15),
    @(103, 44372, 'Especial', 150, 15000, 15000, 15000, '$/caja 15 kilos', 'Región de O''Higgins', 1000, 15),
    @(104, 44372, 'Primera',  300, 13000, 13500, 13250, '$/caja 15 kilos', 'Región de O''Higgins', 883,  15),
    @(105, 44351, 'Especial', 100, 17000, 17000, 17000, '$/caja 15 kilos', 'Región de O''Higgins', 1133, 15),
    @(106, 44351, 'Primera',  400, 14000, 14500, 14250, '$/caja 15 kilos', 'Región de O''Higgins', 950,  15),
    @(107, 44365, 'Primera',  600, 13000, 16000, 14167, '$/caja 15 kilos', 'Región de O''Higgins', 944,  15),
    @(108, 44306, 'Especial', 100, 18000, 18000, 18000, '$/caja 15 kilos', 'Provincia de Cachapoal', 1200, 15),
    @(109, 44306, 'Primera',  400, 14000, 14500, 14250, '$/caja 15 kilos', 'Provincia de Cachapoal', 950,  15),
    @(110, 44411, 'Primera',  350, 14000, 14000, 14000, '$/caja 15 kilos', 'Provincia de Curicó', 933, 15),
    @(111, 44257, 'Primera',  200, 17000, 18000, 17500, '$/caja 15 kilos', 'Región de O''Higgins', 1167, 15),
    @(112, 44376, 'Especial', 150, 15000, 15000, 15000, '$/caja 15 kilos', 'Región de O''Higgins', 1000, 15),
    @(113, 44376, 'Primera',  300, 13000, 13500, 13250, '$/caja 15 kilos', 'Región de O''Higgins', 883,  15),
    @(114, 44292, 'Especial', 200, 18000, 18000, 18000, '$/caja 15 kilos', 'Región de O''Higgins', 1200, 15),
    @(115, 44292, 'Primera',  400, 14000, 15000, 14500, '$/caja 15 kilos', 'Región de O''Higgins', 967,  15),
    @(116, 44358, 'Especial', 200, 16000, 16000, 16000, '$/bandeja 10 kilos', 'Región de O''Higgins', 1600, 10),
    @(117, 44358, 'Primera',  150, 13000, 13000, 13000, '$/bandeja 10 kilos', 'Región de O''Higgins', 1300, 10),
    @(118, 44358, 'Segunda',  150, 13500, 13500, 13500, '$/bandeja 10 kilos', 'Región de O''Higgins', 1350, 10),
    @(119, 44425, 'Especial', 300, 21000, 21000, 21000, '$/caja 15 kilos', 'Provincia de Curicó', 1400, 15),
    @(120, 44425, 'Primera',  150, 14000, 14000, 14000, '$/caja 15 kilos', 'Provincia de Curicó', 933, 15),
    @(121, 44323, 'Especial', 200, 18000, 18000, 18000, '$/caja 15 kilos', 'Región de O''Higgins', 1200, 15),
    @(122, 44323, 'Primera',  400, 14000, 15000, 14500, '$/caja 15 kilos', 'Región de O''Higgins', 967, 15)
)

# Rows 101-119 already exist and only columns D,L,M,N,O,P,Q,R,S,T change.
# Rows 120-122 are brand new and need every column (A-T) populated; they
# reuse the same market/product metadata as the rest of the sheet.
foreach ($entry in $rows) {
    $r = $entry[0]

    if ($r -gt 119) {
        $ws.Cells.Item($r, 1).Value  = 4
        $ws.Cells.Item($r, 2).Value  = 'Feria Lagunitas de Puerto Montt'
        $ws.Cells.Item($r, 3).Value  = 'Los Lagos'
        $ws.Cells.Item($r, 5).Value  = 10
        $ws.Cells.Item($r, 6).Value  = 'Fruta'
        $ws.Cells.Item($r, 7).Value  = 100101
        $ws.Cells.Item($r, 8).Value  = 'Berries'
        $ws.Cells.Item($r, 9).Value  = 100101007
        $ws.Cells.Item($r, 10).Value = 'Kiwi'
        $ws.Cells.Item($r, 11).Value = 'Hayward'
        # New rows need the same date number format used by the rest of
        # column D (copying ".Style" wholesale resets other attributes,
        # so only the number format is transferred explicitly).
        $ws.Cells.Item($r, 4).NumberFormat = $ws.Range("D100").NumberFormat()
    }

    $ws.Cells.Item($r, 4).Value  = $entry[1]
    $ws.Cells.Item($r, 12).Value = $entry[2]
    $ws.Cells.Item($r, 13).Value = $entry[3]
    $ws.Cells.Item($r, 14).Value = $entry[4]
    $ws.Cells.Item($r, 15).Value = $entry[5]
    $ws.Cells.Item($r, 16).Value = $entry[6]
    $ws.Cells.Item($r, 17).Value = $entry[7]
    $ws.Cells.Item($r, 18).Value = $entry[8]
    $ws.Cells.Item($r, 19).Value = $entry[9]
    $ws.Cells.Item($r, 20).Value = $entry[10]
}
